# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker's period window rolls forward by two periods: the account
# statement used to cover periods 2502-2507 (6 periods); it now covers
# 2504-2508 (5 periods). The oldest period's row is dropped and the
# totals are refreshed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Shift each period forward (2507->2504, 2506->2505, 2505->2506, 2504->2507,
# 2503->2508) on rows 16-20 of the detail table.
$ws.Range("E16").Value = "2504"
$ws.Range("E17").Value = "2505"
$ws.Range("E18").Value = "2506"
$ws.Range("E19").Value = "2507"
$ws.Range("E20").Value = "2508"

# Row 21 (period 2502) is being dropped from the table, so row 20 becomes
# the new last row of the table and needs to pick up the closing-border
# formatting that row 21 used to have before that row disappears.
$ws.Range("B21:J21").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows("21:21").Delete()

# Refresh the summary totals: 5 periods now owed at $56,940 each instead
# of the previous 6 periods (5 x 56940 + 17082).
$ws.Range("E11").Value = 284700
$ws.Range("F13").Value = 5
